$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: task renamed to "Construção do pipeline de exemplo"
# and its "X" marker moves from column H to column J
$ws.Range("F6").Value = "Construção do pipeline de exemplo"
$ws.Range("H6").Value = $null
$ws.Range("J6").Value = "X"

# Row 7: new task "Início do pipeline do appraisal" with its "X" marker in column J
$ws.Range("F7").Value = "Início do pipeline do appraisal"
$ws.Range("J7").Value = "X"
